$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 1.09
$ws.Range("J2").Value = 2.52
$ws.Range("M2").Value = 1.02
$ws.Range("O2").Value = 1.08

# Row 3
$ws.Range("F3").Value = 2.34
$ws.Range("G3").Value = 2.56
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 3.6
$ws.Range("J3").Value = 3.15
$ws.Range("K3").Value = 3.4
$ws.Range("L3").Value = 1.49
$ws.Range("N3").Value = 2.68
$ws.Range("P3").Value = 1.59
$ws.Range("Q3").Value = 2.26
$ws.Range("R3").Value = 1.22
$ws.Range("S3").Value = 4.8
$ws.Range("T3").Value = 1.92
$ws.Range("V3").Value = 1.4
$ws.Range("W3").Value = 1.64
$ws.Range("X3").Value = 10.5
$ws.Range("Y3").Value = 11.5
$ws.Range("Z3").Value = 25
$ws.Range("AA3").Value = 80
$ws.Range("AB3").Value = 8.8
$ws.Range("AC3").Value = 7.6
$ws.Range("AD3").Value = 16.5
$ws.Range("AE3").Value = 55
$ws.Range("AF3").Value = 16
$ws.Range("AG3").Value = 12.5
$ws.Range("AH3").Value = 22
$ws.Range("AI3").Value = 90
$ws.Range("AJ3").Value = 38
$ws.Range("AK3").Value = 36
$ws.Range("AM3").Value = 190
$ws.Range("AN3").Value = 36
$ws.Range("AO3").Value = 1000

# Row 4
$ws.Range("F4").Value = 1.76
$ws.Range("G4").Value = 1.89
$ws.Range("H4").Value = 5.1
$ws.Range("I4").Value = 7.2
$ws.Range("J4").Value = 3.1
$ws.Range("K4").Value = 3.65
$ws.Range("L4").Value = 1.4
$ws.Range("M4").Value = 1.09
$ws.Range("N4").Value = 2.98
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 1.66
$ws.Range("Q4").Value = 2.2
$ws.Range("R4").Value = 1.23
$ws.Range("S4").Value = 3.85
$ws.Range("T4").Value = 2.04
$ws.Range("U4").Value = 1.76
$ws.Range("V4").Value = 1.18
$ws.Range("W4").Value = 2.12
$ws.Range("AB4").Value = 29
$ws.Range("AC4").Value = 42
$ws.Range("AG4").Value = 40
$ws.Range("AJ4").Value = 900

# Row 5
$ws.Range("F5").Value = 1.77
$ws.Range("G5").Value = 1.79
$ws.Range("H5").Value = 5.6
$ws.Range("I5").Value = 5.9
$ws.Range("J5").Value = 3.8
$ws.Range("K5").Value = 3.85
$ws.Range("O5").Value = 1.38
$ws.Range("P5").Value = 1.84
$ws.Range("Q5").Value = 2.14
$ws.Range("R5").Value = 1.31
$ws.Range("T5").Value = 2.04
$ws.Range("U5").Value = 1.9
$ws.Range("V5").Value = 1.2
$ws.Range("W5").Value = 2.26
$ws.Range("X5").Value = 12.5
$ws.Range("Y5").Value = 17
$ws.Range("Z5").Value = 42
$ws.Range("AA5").Value = 160
$ws.Range("AC5").Value = 8.4
$ws.Range("AD5").Value = 22
$ws.Range("AE5").Value = 90
$ws.Range("AF5").Value = 9.4
$ws.Range("AH5").Value = 22
$ws.Range("AI5").Value = 95
$ws.Range("AJ5").Value = 17.5
$ws.Range("AK5").Value = 19.5
$ws.Range("AM5").Value = 140
$ws.Range("AN5").Value = 13
$ws.Range("AO5").Value = 120

# Row 6
$ws.Range("F6").Value = 4.8
$ws.Range("G6").Value = 5.8
$ws.Range("H6").Value = 1.76
$ws.Range("I6").Value = 1.86
$ws.Range("J6").Value = 3.55
$ws.Range("K6").Value = 4.3
$ws.Range("L6").Value = 1.33
$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 3.55
$ws.Range("O6").Value = 1.32
$ws.Range("P6").Value = 1.89
$ws.Range("Q6").Value = 1.95
$ws.Range("R6").Value = 1.33
$ws.Range("S6").Value = 3.3
$ws.Range("T6").Value = 1.83
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = 2.16
$ws.Range("W6").Value = 1.22
$ws.Range("Y6").Value = 10.5
$ws.Range("Z6").Value = 13.5
$ws.Range("AC6").Value = 10.5
$ws.Range("AI6").Value = 55
$ws.Range("AK6").Value = 85
$ws.Range("AL6").Value = 110
